$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds "price" values that are stored as plain text (e.g. "60.636.87",
# "0.0558"), not numbers. If assigned directly, Excel reinterprets numeric-looking
# strings as numbers (rounding/precision/format changes). To keep them as text,
# write with a leading apostrophe (forces text entry, like typing it in the UI),
# then restore the "Normal" style so no extra number-format/style is introduced.

$ws.Range('D2').Value = '''60.636.87'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = '  +0.70%  '
$ws.Range('D3').Value = '''2.616.53'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = '  +0.98%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').Value = '''516.10'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  +1.73%  '
$ws.Range('D6').Value = '''154.69'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.47%  '
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('D8').Value = '''0.600'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  +2.67%  '
$ws.Range('D9').Value = '''2.630.52'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  +0.19%  '
$ws.Range('E10').Value = '  +3.74%  '
$ws.Range('E11').Value = '  +0.08%  '
$ws.Range('E12').Value = '  +1.54%  '
$ws.Range('D13').Value = '''0.130'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = '  +1.97%  '
$ws.Range('D14').Value = '''3.073.58'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = '  +0.30%  '
$ws.Range('D15').Value = '''60.637.59'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('D16').Value = '''21.73'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  +0.07%  '
$ws.Range('E17').Value = '  +1.03%  '
$ws.Range('D18').Value = '''2.622.91'
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = '  +0.31%  '
$ws.Range('D19').Value = '''4.75'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = '  -0.58%  '
$ws.Range('D20').Value = '''358.04'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +4.38%  '
$ws.Range('E21').Value = '  +2.69%  '
$ws.Range('E22').Value = '  +1.99%  '
$ws.Range('E23').Value = '  +0.23%  '
$ws.Range('D24').Value = '''61.22'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  +2.18%  '
$ws.Range('E25').Value = '  +1.16%  '
$ws.Range('D26').Value = '''2.734.87'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.86%  '
$ws.Range('E27').Value = '  +0.83%  '
$ws.Range('D28').Value = '''0.997'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +0.41%  '
$ws.Range('E29').Value = '  -1.19%  '
$ws.Range('E30').Value = '  -2.14%  '
$ws.Range('E31').Value = '  +0.22%  '
$ws.Range('D32').Value = '''19.46'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  +0.65%  '
$ws.Range('E33').Value = '  +1.34%  '
$ws.Range('B34').Value = 'Monero'
$ws.Range('C34').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D34').Value = '''152.48'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.28%  '
$ws.Range('B35').Value = 'Aptos'
$ws.Range('C35').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D35').Value = '''5.93'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = '  +3.58%  '
$ws.Range('D36').Value = '''4.05'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  +1.68%  '
$ws.Range('E37').Value = '  -0.81%  '
$ws.Range('D38').Value = '''0.889'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = '  +5.20%  '
$ws.Range('E39').Value = '  +1.08%  '
$ws.Range('D40').Value = '''0.848'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  -0.10%  '
$ws.Range('D41').Value = '''36.23'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +1.81%  '
$ws.Range('E42').Value = '  -1.03%  '
$ws.Range('D43').Value = '''291.55'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -3.28%  '
$ws.Range('D44').Value = '''0.102'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +1.24%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = '''0.0558'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -1.70%  '
$ws.Range('B47').Value = 'FirstDigitalUSD'
$ws.Range('C47').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D47').Value = '''0.996'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.41%  '
$ws.Range('D48').Value = '''19.70'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  -0.43%  '
$ws.Range('E49').Value = '  +0.80%  '
$ws.Range('E50').Value = '  +0.80%  '
$ws.Range('B51').Value = 'InjectiveProtocol'
$ws.Range('C51').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D51').Value = '''19.31'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = '  +11.84%  '
